$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.536975383758545
$ws.Range("B1").Value = 3.669739007949829
$ws.Range("C1").Value = 5.76205587387085
$ws.Range("D1").Value = 1.412290096282959
$ws.Range("E1").Value = 0.8247694969177246
